# "Add files via upload" - append a new work-log row (row 14) to the
# "Eetu Pihamäki" worksheet, describing work done on 1.10.2018.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# Pvm (date) - 1.10.2018
$ws.Range("A14").Value = 43374
# Aloitusklo / Lopetusklo (start/end time of day, stored as day fractions)
$ws.Range("B14").Value = 0.72152777777777777
$ws.Range("C14").Value = 0.79861111111111116
# Sprint column
$ws.Range("E14").Value = 1
# Tehtävä (task) description
$ws.Range("F14").Value = "1h 10 min etsin tiedot 10 järjestelmään, vaatimukseen: `"Tunnusten jäädytys ja poisto`". 30 min järjestelmän valintaa -> löytyy GitHubista `"Vertailutaulukko.xlsx. 5 min esityslistan teko -> lähetetty osallistujille. 5 min työaikakirjanpidon päivitys moodleen. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2028.9.2018%2C%201.10.2018.txt"

# The row grows tall to fit the wrapped description, same as the other
# long task-description rows above it.
$ws.Rows.Item(14).RowHeight = 120

# The active cell in the sheet moves along with the newly entered row.
$ws.Range("F14").Select() | Out-Null

$wb.Save()
